# Normalize the "Recorded By" (column G) values: when the comma-separated
# list of recorders contains an exact "System" entry, move it to the end
# of the list (keeps the other entries, including any differently-cased
# "system", in their original relative order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val.GetType().Name -ne "String") { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "

    $found = $false
    $newParts = @()
    foreach ($p in $parts) {
        if ((-not $found) -and $p.Equals("System")) {
            $found = $true
        } else {
            $newParts += $p
        }
    }

    if ($found) {
        $newParts += "System"
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
